# Añadido borrado múltiple (iss. #4)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new "Tipo" column (C) with "balance" for every existing account row
$ws.Range("C1").Value = "Tipo"
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 3).Value = "balance"
}

# Rename the worksheet from "Cuentas" to "cuentas" (defined name reference
# updates automatically to follow the renamed sheet)
$ws.Name = "cuentas"

# Re-fit the A/B columns now that the layout changed
$ws.Columns("A:B").AutoFit()
